$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.213.36"
$ws.Range("E2").Value = "  -0.39%  "
$ws.Range("D3").Value = "3.050.27"
$ws.Range("E3").Value = "  +1.26%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "515.11"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.33%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.40"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.36%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.440"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.39%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.22"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.54%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.110"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.45%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.378"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.40%  "
$ws.Range("D12").Value = "3.570.47"
$ws.Range("E12").Value = "  +1.21%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.126"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.03%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.08"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.78%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000167"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.26%  "
$ws.Range("D16").Value = "57.165.51"
$ws.Range("E16").Value = "  -0.50%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.18"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.63%  "
$ws.Range("D18").Value = "3.046.61"
$ws.Range("E18").Value = "  +1.24%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.42"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.85%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.16"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.60%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "331.33"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.05%  "
$ws.Range("E22").Value = "  +0.25%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.508"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.85%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.52"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.77%  "
$ws.Range("D25").Value = "3.167.58"
$ws.Range("E25").Value = "  +0.99%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.24%  "
$ws.Range("E27").Value = "  -1.27%  "
$ws.Range("D28").Value = "0.0₃0891"
$ws.Range("E28").Value = "  -2.87%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.76"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.22"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.52%  "
$ws.Range("E31").Value = "  +0.20%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.21"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.69%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.83"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.40%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.72"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.76%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "150.64"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.08%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.97"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.62%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.28"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.32%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "25.34"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.86%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0678"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.31%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.91"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.28%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "36.75"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.84%  "
$ws.Range("E42").Value = "  -0.02%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.665"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.29%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.40"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.75%  "
$ws.Range("D45").Value = "2.203.95"
$ws.Range("E45").Value = "  -0.88%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.12"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.13%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.952"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.86%  "
$ws.Range("E48").Value = "  +1.31%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "20.15"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.22%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.185"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.42%  "
$ws.Range("E51").Value = "  +6.34%  "
